$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.459.59"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.568.59"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D5").Value = "'208.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.78%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'22.19"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.791.08"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.585.62"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  -0.91%  "
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "'63.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("D17").Value = "27.466.05"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'214.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.26%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "'7.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").Value = "'152.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("D28").Value = "'15.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "1.379.25"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'0.954"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.65%  "
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("D43").Value = "'1.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("D44").Value = "'64.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").Value = "1.703.77"
$ws.Range("D48").Value = "'85.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.29%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -0.79%  "
